# Rebuild Sheet1's data table:
#  - drop the old "FICHA" (column G) tracking column entirely
#  - replace the single sample row with five fresh records
#  - move the selection to E17 (matches the author's last cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully remove the old column-G cells (content + formatting) for the header
# and the one data row that used it, so no stray <c> survives there.
$ws.Range("G1:G2").Clear()

# New records: TIPO DE DOCUMENTO, NUMERO DE DOCUMENTO, NOMBRES Y APELLIDOS, DIA, MES, AÑO
$data = @(
    @("CC", 28901342,   "NA", 4,  "NOVIEMBRE", 1978),
    @("CC", 20290328,   "NA", 16, "ENERO",     1963),
    @("CC", 1110545531, "NA", 2,  "ABRIL",     2012),
    @("CC", 24486494,   "NA", 17, "ENERO",     1976),
    @("CC", 94516719,   "NA", 19, "NOVIEMBRE", 1996)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

$ws.Range("E17").Select()
